$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 1-13 (new A values + new B labels) ---
$ws.Range("A1").Value = 100
$ws.Range("B1").Value = "Trial start free water"

$ws.Range("A2").Value = 102
$ws.Range("B2").Value = "Trial: Reward for lick"

$ws.Range("A3").Value = 110
$ws.Range("B3").Value = "Trial aluminum start"

$ws.Range("A4").Value = 111
$ws.Range("B4").Value = "Trial: Reward always"

$ws.Range("A5").Value = 112
$ws.Range("B5").Value = "Trial: Reward for lick"

$ws.Range("A6").Value = 130
$ws.Range("B6").Value = "Trial no Object start"

$ws.Range("A7").Value = 132
$ws.Range("B7").Value = "Trial: CR for no lick"

$ws.Range("A8").Value = 131
$ws.Range("B8").Value = "Trial: No reward"

$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Water no lick"

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Hit"

$ws.Range("A11").Value = 0
$ws.Range("B11").Value = "Miss"

$ws.Range("A12").Value = 66
$ws.Range("B12").Value = "Motor fwd"

$ws.Range("A13").Value = 44
$ws.Range("B13").Value = "motor bck"

# --- New rows 14-19 ---
$ws.Range("A14").Value = 600
$ws.Range("B14").Value = "motor at whisker"

$ws.Range("A15").Value = 400
$ws.Range("B15").Value = "motor at start"

$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Aluminum position"

$ws.Range("A17").Value = 30
$ws.Range("B17").Value = "at non position"

$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "correct rejection"

$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "False alarm"

# --- Column B width ---
$ws.Range("B1").ColumnWidth = 26.15

# --- Selection ---
$ws.Range("B4").Select()
